$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet is a "latest price on top" log. A new circular (20-Sep-2025)
# was published, so a new row is inserted right under the header and
# every existing data row shifts down by one. The newest row gets the
# next Sl.no (34) and the existing rows keep their own Sl.no/price/date
# values (they just move to row+1).
# ------------------------------------------------------------------

# Hyperlinks in this engine are a flat, worksheet-wide collection that
# is NOT re-anchored when rows are inserted/copied, so start by wiping
# them out; they will be rebuilt from scratch, in final position, below.
$ws.Range("A1").Hyperlinks.Delete()

# Insert the new row under the header - pushes rows 2:34 down to 3:35.
$ws.Rows(2).Insert()

# Pick up the same cell styling as the (now) row below (row 3, which is
# the former row 2) so the new row matches the rest of the table.
$ws.Range("A3:F3").Copy($ws.Range("A2:F2"))

# Fill in the brand-new top row with the latest circular's data.
$ws.Range("A2").Value = 34
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 270.5
$ws.Range("E2").Value = "20.09.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-september-2025.pdf"

# Rebuild every hyperlink on column F, now that all rows sit in their
# final position (row 15 - the 12-Aug-2025 circular - gains a link it
# never had before; every other link simply moved down by one row).
$links = @(
    @{ Row = 2;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-september-2025.pdf" },
    @{ Row = 3;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-18-september-2025.pdf" },
    @{ Row = 4;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-17-september-2025.pdf" },
    @{ Row = 5;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-september-2025.pdf" },
    @{ Row = 6;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-september-2025.pdf" },
    @{ Row = 7;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf" },
    @{ Row = 8;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf" },
    @{ Row = 9;  Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf" },
    @{ Row = 10; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf" },
    @{ Row = 11; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf" },
    @{ Row = 12; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf" },
    @{ Row = 13; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf" },
    @{ Row = 14; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf" },
    @{ Row = 15; Url = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf" }
)

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 6)
    $ws.Hyperlinks.Add($cell, $link.Url)
    # Adding a hyperlink swaps in the built-in blue/underlined "Hyperlink"
    # style; the source sheet does not use that style anywhere; restore
    # the plain data-row style (column A in the same row never becomes a
    # hyperlink, so it is a safe, always-on-style format donor).
    $ws.Cells.Item($link.Row, 1).Copy()
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
